$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.09620000000001
$ws.Range("C4").Value = -10.8062
$ws.Range("D4").Value = -6.8777

$ws.Range("C5").Value = -14.66300000000002

$ws.Range("A7").Value = -21.58160000000001

$ws.Range("C8").Value = -11.80969999999999

$ws.Range("D9").Value = -7.870400000000002

$ws.Range("A16").Value = -20.19799999999998
$ws.Range("C16").Value = -12.0892

$ws.Range("D18").Value = -8.421099999999992
